# Append one new data row (row 44) to the daily price sheet, mirroring the
# layout/formatting of the existing rows (e.g. row 43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

$ws.Range("A$row").Value = 1
$ws.Range("B$row").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C$row").Value = "Arica y Parinacota"

# Fecha column: numeric date serial, formatted like the rest of column D.
$ws.Range("D$row").Value = 44448
$ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("E$row").Value = 15
$ws.Range("F$row").Value = 100112038
$ws.Range("G$row").Value = "Cebollín baby"
$ws.Range("H$row").Value = "Sin especificar"
$ws.Range("I$row").Value = "Primera"
$ws.Range("J$row").Value = 270
$ws.Range("K$row").Value = 1900
$ws.Range("L$row").Value = 2000
$ws.Range("M$row").Value = 1950
$ws.Range("N$row").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O$row").Value = "Región de Arica y Parinacota"
$ws.Range("P$row").Value = 975
$ws.Range("Q$row").Value = 2
$ws.Range("R$row").Value = "Hortaliza"
